$d = $word.ActiveDocument

$d.Content.Find.Execute("ESTO ES UNA PRUEBA Y UN MODELO DE EJEMPLO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Prueba 1.", 2)

# Remove every paragraph after the first one (the empty separator
# paragraphs plus the "LOREM IPSUM" paragraphs). Word never lets the
# very last paragraph mark of the document be deleted, so loop until
# only two paragraphs remain (the first one + that trailing mark)...
$target = $d.Paragraphs.Count - 2
$i = 0
while ($i -lt $target) {
    $d.Paragraphs(2).Range.Delete()
    $i = $i + 1
}

# ...then delete paragraph 1's own mark so its text merges into the
# document's final paragraph mark, leaving a single paragraph behind.
$p1 = $d.Paragraphs(1).Range
$markStart = $p1.End - 1
$d.Range($markStart, $p1.End).Delete()
